# Update the cryptos price list with the latest scraped values.
# Values are written with a leading apostrophe (quote-prefix) via .Formula
# so Excel stores them as text (matching the source inlineStr cells)
# instead of auto-converting numeric-looking strings to numbers; the
# follow-up .Style = "Normal" keeps the cell on the default style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'29.375.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Formula = "'  +0.28%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Formula = "'1.884.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Formula = "'  +0.57%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Formula = "'0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Formula = "'  +0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Formula = "'0.7140"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Formula = "'  +0.37%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Formula = "'242.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Formula = "'  -0.06%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Formula = "'0.9997"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Formula = "'0.08072"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Formula = "'  +4.30%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Formula = "'0.3142"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Formula = "'  +1.02%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Formula = "'25.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Formula = "'  +1.25%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Formula = "'0.08367"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Formula = "'  -1.15%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Formula = "'1.875.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Formula = "'  +0.11%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Formula = "'Polygon"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Formula = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Formula = "'0.7233"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Formula = "'  +1.68%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Formula = "'Polkadot"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Formula = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Formula = "'5.257"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Formula = "'  +1.03%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Formula = "'91.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Formula = "'  +0.62%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Formula = "'6.289"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Formula = "'  +5.00%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Formula = "'0.000008466"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Formula = "'  +1.82%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Formula = "'29.359.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Formula = "'  +0.23%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Formula = "'241.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Formula = "'  -0.41%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Formula = "'13.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Formula = "'  +0.49%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Formula = "'2.113.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Formula = "'  -0.37%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Formula = "'1.0000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Formula = "'  +0.04%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Formula = "'7.810"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Formula = "'  -0.07%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Formula = "'1.001"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Formula = "'  +0.05%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Formula = "'0.1593"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Formula = "'  -1.03%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("B26").Formula = "'Monero"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Formula = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Formula = "'163.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Formula = "'  +0.19%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("B27").Formula = "'Cosmos"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Formula = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Formula = "'9.089"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Formula = "'  +0.71%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Formula = "'18.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Formula = "'  +0.49%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Formula = "'1.508"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Formula = "'  -0.39%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Formula = "'4.437"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Formula = "'  +0.66%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Formula = "'4.361"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Formula = "'  +0.99%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Formula = "'  -3.84%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Formula = "'0.05385"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Formula = "'  +2.40%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Formula = "'1.956"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Formula = "'  +1.81%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Formula = "'0.7540"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Formula = "'  +1.32%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Formula = "'1.181"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Formula = "'  +0.69%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Formula = "'  +0.57%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Formula = "'0.01883"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Formula = "'  +1.31%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Formula = "'1.281.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Formula = "'  +9.52%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Formula = "'  +1.15%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Formula = "'6.581"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Formula = "'  +3.54%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Formula = "'110.55"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Formula = "'  +3.67%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Formula = "'73.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Formula = "'  +0.73%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Formula = "'0.8920"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Formula = "'  +0.82%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Formula = "'1.000"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Formula = "'  +0.04%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Formula = "'  +7.56%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Formula = "'2.004.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Formula = "'  -0.75%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Formula = "'1.808"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Formula = "'  +0.00%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Formula = "'0.5218"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Formula = "'  +0.27%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Formula = "'9.495"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Formula = "'  +1.28%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Formula = "'0.4378"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Formula = "'  +1.80%  "
$ws.Range("E51").Style = "Normal"
